$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("datasets")

# Preserve the existing number/cell formatting (quote-prefix style) of I2:I52
# before we overwrite the formulas, then restore it afterwards via a
# formats-only paste (xlPasteFormats = -4122) through a scratch column.
$ws.Range("I2:I52").Copy() | Out-Null
$ws.Range("Z2:Z52").PasteSpecial(-4122) | Out-Null

# Update the dataset path prefix from "~/VREFolders/..." to
# "~/workspace/VREFolders/..." in every formula of column I (rows 2-52).
for ($r = 2; $r -le 52; $r++) {
    $ws.Range("I$r").Formula = '=+_xlfn.CONCAT("~/workspace/VREFolders/ITINERIS_EV/DATI/",H' + $r + ')'
}

# Restore the original cell formatting that got reset by the formula rewrite.
$ws.Range("Z2:Z52").Copy() | Out-Null
$ws.Range("I2:I52").PasteSpecial(-4122) | Out-Null

# Clean up the scratch column used to stash the formatting.
$ws.Range("Z2:Z52").Clear() | Out-Null

# Match the final selection left behind in the sheet.
$ws.Activate() | Out-Null
$ws.Range("I2:I52").Select() | Out-Null
